$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column D header cell (D2, blank but styled like the header row)
# and D3 ("Note" header) both use the same bold header style as B3/C3.
$ws.Range("B3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)

# Header row additions / edits
$ws.Range("D3").Value = "Note"

# AFDD-1 section
$ws.Range("D6").Value = "The passive diagnostic does not deterimine which temperature sensor is faulty."
$ws.Range("D11").Value = "Missing supply-fan status, outdoor-air temperature, mixed-air temperature, or return-air temperature sensor readings."

# AFDD-2 section (title text corrected)
$ws.Range("C14").Value = "AFDD-2 (Unit Not Economizing When it Should)"
$ws.Range("D17").Value = "Unit is not fully utilizing economzing or a mechanical issue is causing outdoor-air fraction to be too low."
$ws.Range("D19").Value = "Mixed-air temperature (MAT) is not trended and the heating or cooling is active.  Since discharge-air temperature (DAT) is used instead of MAT the heating and cooling cannot be active (unfavorable condition)."
$ws.Range("D20").Value = "Missing supply-fan status, outdoor-air tempeature, return-air temperature, mixed-air temperature, outdoor damper command, or compressor command."
$ws.Range("D21").Value = "OAF calculation is not in the range of [0,1.25]"

# AFDD-3 section (title text corrected)
$ws.Range("C24").Value = "AFDD-3 (Unit Economizing When it Should Not)"
$ws.Range("D27").Value = "Since some units utilize integrated economizer logic and only the compressor command is known a unit maybe in a cooling mode and economizing when the compressor is off."
$ws.Range("D28").Value = "Missing outdoor-air tempeature, return-air temperature, mixed-air temperature, outdoor damper command, or compressor command."

# AFDD-4 section
$ws.Range("D36").Value = "Since some units utilize integrated economizer logic and only the compressor command is known a unit maybe in a cooling mode and economizing when the compressor is off."
$ws.Range("D37").Value = "Missing supply-fan status, outdoor-air tempeature, return-air temperature, mixed-air temperature, outdoor damper command, or compressor command."
$ws.Range("D38").Value = "OAF calculation is not in the range of [0,1.25]"

# AFDD-5 section
$ws.Range("D44").Value = "The difference of the outdoor-air temperature and mixed-air temperature should be at least 5F for OAF calculation."
$ws.Range("C45").Value = "Damper is not at minimum when is  not be (Fault)"
$ws.Range("D45").Value = "Outdoor-air damper is significantly above the minimum (correct for current conditions) command."
$ws.Range("D46").Value = "Since some units utilize integrated economizer logic and only the compressor command is known a unit maybe in a cooling mode and economizing when the compressor is off."
$ws.Range("D47").Value = "Missing supply-fan status, outdoor-air tempeature, return-air temperature, mixed-air temperature, outdoor damper command, or compressor command."
$ws.Range("C48").Value = "Damper is at minimum but OAF calculation led to  an unexpected value (No Fault)"
$ws.Range("D48").Value = "OAF calculation is not in the range of [0,1.25]"

# Final selection / view state left by the editor
$ws.Range("E39").Select()
